# Update the build timestamp embedded in the version strings from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet ---
$a2 = $aboutSheet.Range("A2").Value2
$aboutSheet.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $aboutSheet.Range("A6").Value2
$aboutSheet.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- Boundaries and methane sources sheet (column S, rows 2-7) ---
foreach ($row in 2..7) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value2
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
